$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2847.1667
$ws.Range("I64").Value = 2816.6
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 2816.6
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -2568.6
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 2847.1667
$ws.Range("I67").Value = 2816.6
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 2816.6
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -1958.6
$ws.Range("N67").Value = -4716
$ws.Range("H70").Value = 32975.25
$ws.Range("I70").Value = 950.5
$ws.Range("J70").Value = 65000
$ws.Range("K70").Value = 2851.5
$ws.Range("L70").Value = 195000
$ws.Range("M70").Value = -2581.5
$ws.Range("N70").Value = -195540
$ws.Range("H73").Value = 32975.25
$ws.Range("I73").Value = 950.5
$ws.Range("J73").Value = 65000
$ws.Range("K73").Value = 2851.5
$ws.Range("L73").Value = 195000
$ws.Range("M73").Value = -1915.5
$ws.Range("N73").Value = -196872
$ws.Range("H94").Value = 2552.7778
$ws.Range("I94").Value = 2552.7778
$ws.Range("K94").Value = 2552.7778
$ws.Range("M94").Value = -2101.7778
$ws.Range("H106").Value = 2188.3333
$ws.Range("J106").Value = 1231.3334
$ws.Range("L106").Value = 1231.3334
$ws.Range("N106").Value = -2493.3334
$ws.Range("H112").Value = 1486.6285
$ws.Range("J112").Value = 1500.9412
$ws.Range("L112").Value = 4502.8236
$ws.Range("N112").Value = -6718.8236
$ws.Range("H135").Value = 516.06665
$ws.Range("I135").Value = 170.6
$ws.Range("J135").Value = 1207
$ws.Range("K135").Value = 1535.4
$ws.Range("L135").Value = 10863
$ws.Range("M135").Value = 999.6000000000001
$ws.Range("N135").Value = -15933
$ws.Range("H138").Value = 2534.9558
$ws.Range("J138").Value = 2031.4791
$ws.Range("L138").Value = 6094.4373
$ws.Range("N138").Value = -16374.4373

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4220.8667
$ws.Range("I32").Value = 3360.0544
$ws.Range("K32").Value = 3360.0544
$ws.Range("M32").Value = -3073.0544
$ws.Range("H45").Value = 1498.1177
$ws.Range("I45").Value = 996
$ws.Range("K45").Value = 996
$ws.Range("M45").Value = -619
$ws.Range("H61").Value = 7228.8887
$ws.Range("I61").Value = 8055.5386
$ws.Range("K61").Value = 8055.5386
$ws.Range("M61").Value = -7843.5386
$ws.Range("H74").Value = 1270.381
$ws.Range("J74").Value = 4594.75
$ws.Range("L74").Value = 4594.75
$ws.Range("N74").Value = -6342.75
$ws.Range("H77").Value = 1270.381
$ws.Range("J77").Value = 4594.75
$ws.Range("L77").Value = 22973.75
$ws.Range("N77").Value = -31709.75
$ws.Range("H88").Value = 3727.1667
$ws.Range("I88").Value = 2500
$ws.Range("K88").Value = 2500
$ws.Range("M88").Value = -2094
$ws.Range("H91").Value = 3727.1667
$ws.Range("I91").Value = 2500
$ws.Range("K91").Value = 2500
$ws.Range("M91").Value = -1096
$ws.Range("H109").Value = 58971.25
$ws.Range("J109").Value = 58971.25
$ws.Range("L109").Value = 58971.25
$ws.Range("N109").Value = -61745.25
$ws.Range("H132").Value = 3013.7
$ws.Range("I132").Value = 2465.8333
$ws.Range("J132").Value = 3835.5
$ws.Range("K132").Value = 7397.499899999999
$ws.Range("L132").Value = 11506.5
$ws.Range("M132").Value = -4867.499899999999
$ws.Range("N132").Value = -16566.5
$ws.Range("H136").Value = 7228.8887
$ws.Range("I136").Value = 8055.5386
$ws.Range("K136").Value = 24166.6158
$ws.Range("M136").Value = -21616.6158

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 86166.125
$ws.Range("I86").Value = 3429.875
$ws.Range("K86").Value = 3429.875
$ws.Range("M86").Value = -2306.875
$ws.Range("H89").Value = 86166.125
$ws.Range("I89").Value = 3429.875
$ws.Range("K89").Value = 17149.375
$ws.Range("M89").Value = -11533.375
$ws.Range("H134").Value = 5487.0386
$ws.Range("I134").Value = 5770.696
$ws.Range("K134").Value = 17312.088
$ws.Range("M134").Value = -14777.088

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2304.2856
$ws.Range("I31").Value = 1899.2
$ws.Range("J31").Value = 2529.3333
$ws.Range("K31").Value = 1899.2
$ws.Range("L31").Value = 2529.3333
$ws.Range("M31").Value = -1604.2
$ws.Range("N31").Value = -3119.3333
$ws.Range("H34").Value = 2304.2856
$ws.Range("I34").Value = 1899.2
$ws.Range("J34").Value = 2529.3333
$ws.Range("K34").Value = 1899.2
$ws.Range("L34").Value = 2529.3333
$ws.Range("M34").Value = -1697.2
$ws.Range("N34").Value = -2933.3333
$ws.Range("H62").Value = 3700
$ws.Range("I62").Value = 3650
$ws.Range("J62").Value = 3800
$ws.Range("K62").Value = 3650
$ws.Range("L62").Value = 3800
$ws.Range("M62").Value = -3026
$ws.Range("N62").Value = -5048
$ws.Range("H65").Value = 3700
$ws.Range("I65").Value = 3650
$ws.Range("J65").Value = 3800
$ws.Range("K65").Value = 18250
$ws.Range("L65").Value = 19000
$ws.Range("M65").Value = -15130
$ws.Range("N65").Value = -25240
$ws.Range("H86").Value = 3165.1428
$ws.Range("I86").Value = 3177.8462
$ws.Range("K86").Value = 3177.8462
$ws.Range("M86").Value = -2054.8462
$ws.Range("H89").Value = 3165.1428
$ws.Range("I89").Value = 3177.8462
$ws.Range("K89").Value = 15889.231
$ws.Range("M89").Value = -10273.231
$ws.Range("H92").Value = 42498.5
$ws.Range("J92").Value = 42498.5
$ws.Range("L92").Value = 42498.5
$ws.Range("N92").Value = -47490.5
$ws.Range("H141").Value = 62247.75
$ws.Range("J141").Value = 58997
$ws.Range("L141").Value = 58997
$ws.Range("N141").Value = -69357

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3302.3914
$ws.Range("I134").Value = 2195.111
$ws.Range("J134").Value = 4014.2144
$ws.Range("K134").Value = 6585.333
$ws.Range("L134").Value = 12042.6432
$ws.Range("M134").Value = -1515.333
$ws.Range("N134").Value = -22182.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1788.9
$ws.Range("I68").Value = 1432.1111
$ws.Range("K68").Value = 1432.1111
$ws.Range("M68").Value = -683.1111000000001
$ws.Range("H71").Value = 1788.9
$ws.Range("I71").Value = 1432.1111
$ws.Range("K71").Value = 7160.5555
$ws.Range("M71").Value = -3416.5555
$ws.Range("H122").Value = 8696.286
$ws.Range("I122").Value = 8920.1
$ws.Range("K122").Value = 26760.3
$ws.Range("M122").Value = -24310.3
$ws.Range("H132").Value = 3313.55
$ws.Range("I132").Value = 1999.8572
$ws.Range("K132").Value = 5999.571599999999
$ws.Range("M132").Value = -3469.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 4798.75
$ws.Range("I45").Value = 3569
$ws.Range("J45").Value = 5208.6665
$ws.Range("K45").Value = 3569
$ws.Range("L45").Value = 5208.6665
$ws.Range("N45").Value = -6190.6665
$ws.Range("M45").Value = -3078
$ws.Range("H122").Value = 20428.635
$ws.Range("I122").Value = 31244.73
$ws.Range("J122").Value = 1680.7333
$ws.Range("K122").Value = 93734.19
$ws.Range("L122").Value = 5042.199900000001
$ws.Range("M122").Value = -91284.19
$ws.Range("N122").Value = -9942.1999
$ws.Range("H132").Value = 2522.0908
$ws.Range("I132").Value = 1615.6666
$ws.Range("J132").Value = 4464.4287
$ws.Range("K132").Value = 4846.9998
$ws.Range("L132").Value = 13393.2861
$ws.Range("M132").Value = -2316.9998
$ws.Range("N132").Value = -18453.2861
